$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Paragraph 1 currently reads:
#   <w:p>
#     <w:pPr>
#       <w:spacing w:line="240" w:lineRule="auto"/>
#       <w:ind w:firstLine="720"/>
#       <w:rPr>...</w:rPr>
#     </w:pPr>
#     <w:r><w:rPr>...</w:rPr><w:t>This is method.</w:t></w:r>
#   </w:p>
#
# Target:
#   - the first-line indent is removed from the paragraph properties
#   - the run text becomes "This is a place for method." split across three
#     runs ("This is ", "a place for ", "method.") with the (reserved,
#     single-instance) "_GoBack" bookmark sitting right before "method."
# --------------------------------------------------------------------------

$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range

# Step 1: rewrite paragraph 1's XML, dropping <w:ind .../> and updating the
# text to the new combined sentence. Paragraph-level ids / rsid attributes
# are carried over unchanged so nothing else about the paragraph shifts.
$p1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
         'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
         'w14:paraId="2F0AE327" w14:textId="4C95BC50" w:rsidR="00BA652F" ' +
         'w:rsidRDefault="00894BEB" w:rsidP="004342C5">' +
           '<w:pPr>' +
             '<w:spacing w:line="240" w:lineRule="auto"/>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>' +
               '<w:color w:val="000000" w:themeColor="text1"/>' +
             '</w:rPr>' +
           '</w:pPr>' +
           '<w:r>' +
             '<w:rPr>' +
               '<w:rFonts w:ascii="Cambria" w:hAnsi="Cambria"/>' +
               '<w:color w:val="000000" w:themeColor="text1"/>' +
             '</w:rPr>' +
             '<w:t>This is a place for method.</w:t>' +
           '</w:r>' +
         '</w:p>'
$p1Range.InsertXML($p1Xml) | Out-Null

# Step 2: re-find the "a place for " span inside paragraph 1 and drop the
# "_GoBack" bookmark right after it. A document only ever carries one
# "_GoBack" bookmark, so adding it here both splits the run in two (the
# bookmark sits between "a place for " and "method.") and implicitly
# removes it from wherever it used to live (the second, empty paragraph).
$placeRange = $d.Content
$placeRange.Find.Execute("a place for ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($placeRange.End, $placeRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Step 3: split "This is " from "a place for " into separate runs. Adding a
# bookmark at the boundary forces the run break; deleting that scratch
# bookmark immediately afterwards leaves the split in place without adding
# any visible bookmark markup.
$thisIsRange = $d.Content
$thisIsRange.Find.Execute("This is ", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0) | Out-Null
$splitRange = $d.Range($thisIsRange.End, $thisIsRange.End)
$d.Bookmarks.Add("ZZSplitTemp", $splitRange)
$d.Bookmarks("ZZSplitTemp").Delete()
